# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The "Rules" sheet row 11 (R40) gets its rule-name cell (B11) replaced
# with the literal text "1" (still a plain text cell, not a number).
#
# A straight $cell.Value = "1" assignment gets auto-coerced to the
# number 1 by Excel's usual "smart" cell-entry logic (since "1" looks
# numeric) - that's correct Excel behavior for typed/assigned input,
# but it's not what we want here: we need literal text "1".
#
# To force literal text without perturbing the cell's existing style,
# compute the text "1" with TEXT(), then convert that formula result to
# a static value via Copy / Paste-Special-Values (xlPasteValues), which
# preserves the String type of the pasted value instead of re-parsing
# it as typed input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")

# Produce the text "1" as a formula result...
$cell.Formula = "=TEXT(1,""0"")"

# ...then freeze it into a literal text value (Paste Special - Values),
# which keeps it text instead of renumbering it.
$cell.Copy()
$cell.PasteSpecial(-4163)
$excel.CutCopyMode = 0
